# Insert a new data row before the current row 245 (shifts old rows
# 245..322 down to 246..323) and populate the new row with a fresh
# "Acelga" price observation for Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 245 downward by one.
$ws.Rows(245).Insert()

# Fill in the newly inserted row 245 with the new observation.
$ws.Cells.Item(245, 1).Value  = 4
$ws.Cells.Item(245, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(245, 3).Value  = "Los Lagos"
$ws.Cells.Item(245, 4).Value  = 45093
$ws.Cells.Item(245, 5).Value  = 10
$ws.Cells.Item(245, 6).Value  = 100112009
$ws.Cells.Item(245, 7).Value  = "Acelga"
$ws.Cells.Item(245, 8).Value  = "Sin especificar"
$ws.Cells.Item(245, 9).Value  = "Primera"
$ws.Cells.Item(245, 10).Value = 95
$ws.Cells.Item(245, 11).Value = 9000
$ws.Cells.Item(245, 12).Value = 9000
$ws.Cells.Item(245, 13).Value = 9000
$ws.Cells.Item(245, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(245, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(245, 16).Value = 750
$ws.Cells.Item(245, 17).Value = 12
$ws.Cells.Item(245, 18).Value = "Hortaliza"
